# Natmi following Dr Hou advice:
# Add the "ECs" sending/target cluster to the Fgf1-Fgfr3 LR-pairs sheet,
# recomputing the full 3x3 (ECs / FAPs / sCs) cluster-pair cross table
# in rows 2-10 (previously only the 2x2 FAPs/sCs cross table existed).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Fgf1"
$ws.Cells.Item(2, 3).Value = "Fgfr3"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 2
$ws.Cells.Item(2, 6).Value = 0.6666666666666666
$ws.Cells.Item(2, 7).Value = 2.004760666666666
$ws.Cells.Item(2, 8).Value = 6.014282
$ws.Cells.Item(2, 9).Value = 0.1200698528618338
$ws.Cells.Item(2, 10).Value = 0.1200698528618338
$ws.Cells.Item(2, 11).Value = 2
$ws.Cells.Item(2, 12).Value = 0.6666666666666666
$ws.Cells.Item(2, 13).Value = 4.220261333333333
$ws.Cells.Item(2, 14).Value = 12.660784
$ws.Cells.Item(2, 15).Value = 0.6739259863235564
$ws.Cells.Item(2, 16).Value = 0.6739259863235564
$ws.Cells.Item(2, 17).Value = 8.460613924120887
$ws.Cells.Item(2, 18).Value = 76.14552531708799
$ws.Cells.Item(2, 19).Value = 0.08091819401763561
$ws.Cells.Item(2, 20).Value = 0.08091819401763561

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Fgf1"
$ws.Cells.Item(3, 3).Value = "Fgfr3"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 2
$ws.Cells.Item(3, 6).Value = 0.6666666666666666
$ws.Cells.Item(3, 7).Value = 2.004760666666666
$ws.Cells.Item(3, 8).Value = 6.014282
$ws.Cells.Item(3, 9).Value = 0.1200698528618338
$ws.Cells.Item(3, 10).Value = 0.1200698528618338
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 1.296447666666667
$ws.Cells.Item(3, 14).Value = 3.889343
$ws.Cells.Item(3, 15).Value = 0.2070274097896007
$ws.Cells.Item(3, 16).Value = 0.2070274097896007
$ws.Cells.Item(3, 17).Value = 2.599067288525111
$ws.Cells.Item(3, 18).Value = 23.391605596726
$ws.Cells.Item(3, 19).Value = 0.02485775063180392
$ws.Cells.Item(3, 20).Value = 0.02485775063180392

# Row 4
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Fgf1"
$ws.Cells.Item(4, 3).Value = "Fgfr3"
$ws.Cells.Item(4, 4).Value = "sCs"
$ws.Cells.Item(4, 5).Value = 2
$ws.Cells.Item(4, 6).Value = 0.6666666666666666
$ws.Cells.Item(4, 7).Value = 2.004760666666666
$ws.Cells.Item(4, 8).Value = 6.014282
$ws.Cells.Item(4, 9).Value = 0.1200698528618338
$ws.Cells.Item(4, 10).Value = 0.1200698528618338
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 0.745494
$ws.Cells.Item(4, 14).Value = 2.236482
$ws.Cells.Item(4, 15).Value = 0.119046603886843
$ws.Cells.Item(4, 16).Value = 0.119046603886843
$ws.Cells.Item(4, 17).Value = 1.494537048436
$ws.Cells.Item(4, 18).Value = 13.450833435924
$ws.Cells.Item(4, 19).Value = 0.01429390821239425
$ws.Cells.Item(4, 20).Value = 0.01429390821239425

# Row 5
$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Fgf1"
$ws.Cells.Item(5, 3).Value = "Fgfr3"
$ws.Cells.Item(5, 4).Value = "ECs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 9.409654999999999
$ws.Cells.Item(5, 8).Value = 28.228965
$ws.Cells.Item(5, 9).Value = 0.5635664696121425
$ws.Cells.Item(5, 10).Value = 0.5635664696121424
$ws.Cells.Item(5, 11).Value = 2
$ws.Cells.Item(5, 12).Value = 0.6666666666666666
$ws.Cells.Item(5, 13).Value = 4.220261333333333
$ws.Cells.Item(5, 14).Value = 12.660784
$ws.Cells.Item(5, 15).Value = 0.6739259863235564
$ws.Cells.Item(5, 16).Value = 0.6739259863235564
$ws.Cells.Item(5, 17).Value = 39.71120315650666
$ws.Cells.Item(5, 18).Value = 357.40082840856
$ws.Cells.Item(5, 19).Value = 0.3798020888922477
$ws.Cells.Item(5, 20).Value = 0.3798020888922476

# Row 6
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Fgf1"
$ws.Cells.Item(6, 3).Value = "Fgfr3"
$ws.Cells.Item(6, 4).Value = "FAPs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 9.409654999999999
$ws.Cells.Item(6, 8).Value = 28.228965
$ws.Cells.Item(6, 9).Value = 0.5635664696121425
$ws.Cells.Item(6, 10).Value = 0.5635664696121424
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 1.296447666666667
$ws.Cells.Item(6, 14).Value = 3.889343
$ws.Cells.Item(6, 15).Value = 0.2070274097896007
$ws.Cells.Item(6, 16).Value = 0.2070274097896007
$ws.Cells.Item(6, 17).Value = 12.19912526888833
$ws.Cells.Item(6, 18).Value = 109.792127419995
$ws.Cells.Item(6, 19).Value = 0.1166737064480716
$ws.Cells.Item(6, 20).Value = 0.1166737064480716

# Row 7
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Fgf1"
$ws.Cells.Item(7, 3).Value = "Fgfr3"
$ws.Cells.Item(7, 4).Value = "sCs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 9.409654999999999
$ws.Cells.Item(7, 8).Value = 28.228965
$ws.Cells.Item(7, 9).Value = 0.5635664696121425
$ws.Cells.Item(7, 10).Value = 0.5635664696121424
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 0.745494
$ws.Cells.Item(7, 14).Value = 2.236482
$ws.Cells.Item(7, 15).Value = 0.119046603886843
$ws.Cells.Item(7, 16).Value = 0.119046603886843
$ws.Cells.Item(7, 17).Value = 7.014841344569999
$ws.Cells.Item(7, 18).Value = 63.13357210113
$ws.Cells.Item(7, 19).Value = 0.06709067427182329
$ws.Cells.Item(7, 20).Value = 0.06709067427182328

# Row 8
$ws.Cells.Item(8, 1).Value = "sCs"
$ws.Cells.Item(8, 2).Value = "Fgf1"
$ws.Cells.Item(8, 3).Value = "Fgfr3"
$ws.Cells.Item(8, 4).Value = "ECs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 5.282203999999999
$ws.Cells.Item(8, 8).Value = 15.846612
$ws.Cells.Item(8, 9).Value = 0.3163636775260238
$ws.Cells.Item(8, 10).Value = 0.3163636775260238
$ws.Cells.Item(8, 11).Value = 2
$ws.Cells.Item(8, 12).Value = 0.6666666666666666
$ws.Cells.Item(8, 13).Value = 4.220261333333333
$ws.Cells.Item(8, 14).Value = 12.660784
$ws.Cells.Item(8, 15).Value = 0.6739259863235564
$ws.Cells.Item(8, 16).Value = 0.6739259863235564
$ws.Cells.Item(8, 17).Value = 22.29228129597866
$ws.Cells.Item(8, 18).Value = 200.630531663808
$ws.Cells.Item(8, 19).Value = 0.2132057034136731
$ws.Cells.Item(8, 20).Value = 0.2132057034136731

# Row 9
$ws.Cells.Item(9, 1).Value = "sCs"
$ws.Cells.Item(9, 2).Value = "Fgf1"
$ws.Cells.Item(9, 3).Value = "Fgfr3"
$ws.Cells.Item(9, 4).Value = "FAPs"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 5.282203999999999
$ws.Cells.Item(9, 8).Value = 15.846612
$ws.Cells.Item(9, 9).Value = 0.3163636775260238
$ws.Cells.Item(9, 10).Value = 0.3163636775260238
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 1.296447666666667
$ws.Cells.Item(9, 14).Value = 3.889343
$ws.Cells.Item(9, 15).Value = 0.2070274097896007
$ws.Cells.Item(9, 16).Value = 0.2070274097896007
$ws.Cells.Item(9, 17).Value = 6.848101050657332
$ws.Cells.Item(9, 18).Value = 61.632909455916
$ws.Cells.Item(9, 19).Value = 0.0654959527097252
$ws.Cells.Item(9, 20).Value = 0.0654959527097252

# Row 10
$ws.Cells.Item(10, 1).Value = "sCs"
$ws.Cells.Item(10, 2).Value = "Fgf1"
$ws.Cells.Item(10, 3).Value = "Fgfr3"
$ws.Cells.Item(10, 4).Value = "sCs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 5.282203999999999
$ws.Cells.Item(10, 8).Value = 15.846612
$ws.Cells.Item(10, 9).Value = 0.3163636775260238
$ws.Cells.Item(10, 10).Value = 0.3163636775260238
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 0.745494
$ws.Cells.Item(10, 14).Value = 2.236482
$ws.Cells.Item(10, 15).Value = 0.119046603886843
$ws.Cells.Item(10, 16).Value = 0.119046603886843
$ws.Cells.Item(10, 17).Value = 3.937851388775999
$ws.Cells.Item(10, 18).Value = 35.440662498984
$ws.Cells.Item(10, 19).Value = 0.0376620214026255
$ws.Cells.Item(10, 20).Value = 0.0376620214026255
